$wb = $excel.ActiveWorkbook

# ---- Sheet7 (data sheet, first tab) ----
$ws = $wb.Worksheets.Item("Sheet7")
$ws.Activate()

# New "TOTALVOLUME KA PENUMPANG" (column D) figures
$dValues = @{
    2 = 16268;  3 = 12058;  4 = 11683;  5 = 5862;   6 = 36315;
    7 = 23762;  8 = 85099;  9 = 25783;  10 = 95410;
    11 = 14966; 12 = 16209; 13 = 58931; 14 = 71342; 15 = 21643;
    16 = 17495; 17 = 33485; 18 = 29607; 19 = 38696; 20 = 52227;
    21 = 36873; 22 = 37292; 23 = 70891; 24 = 129471; 25 = 27540;
    26 = 72036; 27 = 242;   28 = 752;   29 = 6897;  30 = 391;
    31 = 15087; 32 = 8834;  33 = 457;   34 = 135;   35 = 190;
    36 = 482;   37 = 5350
}

foreach ($row in $dValues.Keys) {
    $cell = $ws.Cells.Item($row, 4)
    $cell.Value2 = $dValues[$row]
    $cell.HorizontalAlignment = -4108
    $cell.WrapText = $true
}

# Widen column D (previously auto/best-fit at ~8.78 chars) and let the header
# row shrink back down from its old wrapped height of 66 to 26.4
$ws.Columns("D").ColumnWidth = 20.6
$ws.Rows(1).RowHeight = 26.4

# Update the view: scrolled down with the newly-reviewed block selected
$ws.Range("A27:E29").Select()
$excel.ActiveWindow.ScrollRow = 16

# ---- Sheet3 (reference sheet) ----
$ws2 = $wb.Worksheets.Item("Sheet3")
$ws2.Activate()
$ws2.Range("A1:E40").Select()
$excel.ActiveWindow.ScrollRow = 21
